$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Inscritos (E5) 25 -> 26
$ws.Range("E5").Value = 26

# Row 6: Pagos (F6) 24 -> 25, Inscricoes homologadas (H6) 24 -> 25
$ws.Range("F6").Value = 25
$ws.Range("H6").Value = 25

# Row 10: Inscritos (E10) 22 -> 23, Pagos (F10) 9 -> 10, Inscricoes homologadas (H10) 9 -> 10
$ws.Range("E10").Value = 23
$ws.Range("F10").Value = 10
$ws.Range("H10").Value = 10

# Row 16: Inscritos (E16) 292 -> 294
$ws.Range("E16").Value = 294

# Row 17: Inscritos (E17) 18 -> 19
$ws.Range("E17").Value = 19
